$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.719.90"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "'1.534.02"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'205.84"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  -1.04%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'21.37"
$ws.Range("E8").Value = "  -2.68%  "
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").Value = "'1.750.27"
$ws.Range("E12").Value = "  -1.56%  "
$ws.Range("D13").Value = "'1.532.42"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").Value = "'3.68"
$ws.Range("D15").Value = "'0.507"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "'61.30"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "'26.703.64"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "'212.53"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "'0.0₃0682"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("D20").Value = "'7.21"
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("D23").Value = "'9.10"
$ws.Range("E23").Value = "  -2.64%  "
$ws.Range("E24").Value = "  -3.56%  "
$ws.Range("D25").Value = "'151.97"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Value = "'6.55"
$ws.Range("E26").Value = "  -3.31%  "
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("E30").Value = "  -0.99%  "
$ws.Range("D31").Value = "'0.0453"
$ws.Range("E31").Value = "  -1.92%  "
$ws.Range("D32").Value = "'3.24"
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("D33").Value = "'1.356.82"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("E35").Value = "  -3.40%  "
$ws.Range("D36").Value = "'0.940"
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'0.797"
$ws.Range("E40").Value = "  -1.69%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "'5.68"
$ws.Range("E42").Value = "  +5.42%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'1.74"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'62.61"
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("D47").Value = "'1.664.57"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("D48").Value = "'85.43"
$ws.Range("D49").Value = "'0.0507"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("D50").Value = "'0.0944"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.06%  "
